$wb = $excel.ActiveWorkbook

# --- Productdata sheet ---
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 0.052
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0.04960000000000001
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0.0192
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 0.008
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 0.0144
$ws.Range("E7").Value = 0.026
$ws.Range("E8").Value = 0.0508
$ws.Range("E9").Value = 0.0248

# --- Capacity sheet ---
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 40
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 80
$ws.Range("B7").Value = 15
$ws.Range("B8").Value = 10
$ws.Range("B9").Value = 10

# --- ProcessingTime sheet ---
$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("D4").Value = 2
$ws.Range("E5").Value = 4
$ws.Range("F6").Value = 4
$ws.Range("H8").Value = 2
$ws.Range("I9").Value = 2
